$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 139.633461
$ws.Range("H2").Value = 418.900383
$ws.Range("I2").Value = 0.1723688299390553
$ws.Range("J2").Value = 0.1847819285778738
$ws.Range("M2").Value = 38.10639333333333
$ws.Range("N2").Value = 114.31918
$ws.Range("O2").Value = 0.3831479157160237
$ws.Range("P2").Value = 0.4159903984418967
$ws.Range("Q2").Value = 5320.92758736066
$ws.Range("R2").Value = 47888.34828624594
$ws.Range("S2").Value = 0.0660427579255588
$ws.Range("T2").Value = 0.0768675080939718
$ws.Range("G3").Value = 139.633461
$ws.Range("H3").Value = 418.900383
$ws.Range("I3").Value = 0.1723688299390553
$ws.Range("J3").Value = 0.1847819285778738
$ws.Range("O3").Value = 0.09199521176963764
$ws.Range("P3").Value = 0.09988081163714851
$ws.Range("Q3").Value = 1277.574117284124
$ws.Range("R3").Value = 11498.16705555712
$ws.Range("S3").Value = 0.01585710701272805
$ws.Range("T3").Value = 0.01845616900223564
$ws.Range("G4").Value = 139.633461
$ws.Range("H4").Value = 418.900383
$ws.Range("I4").Value = 0.1723688299390553
$ws.Range("J4").Value = 0.1847819285778738
$ws.Range("M4").Value = 15.023598
$ws.Range("N4").Value = 45.070794
$ws.Range("O4").Value = 0.1510575983904562
$ws.Range("P4").Value = 0.1640058785774412
$ws.Range("Q4").Value = 2097.796985412678
$ws.Range("R4").Value = 18880.1728687141
$ws.Range("S4").Value = 0.02603762148796666
$ws.Range("T4").Value = 0.03030532254164818
$ws.Range("G5").Value = 139.633461
$ws.Range("H5").Value = 418.900383
$ws.Range("I5").Value = 0.1723688299390553
$ws.Range("J5").Value = 0.1847819285778738
$ws.Range("M5").Value = 23.556204
$ws.Range("N5").Value = 47.112408
$ws.Range("O5").Value = 0.236850294013169
$ws.Range("P5").Value = 0.1714350065796238
$ws.Range("Q5").Value = 3289.234292542044
$ws.Range("R5").Value = 19735.40575525227
$ws.Range("S5").Value = 0.04082560804977118
$ws.Range("T5").Value = 0.03167809114154337
$ws.Range("G6").Value = 139.633461
$ws.Range("H6").Value = 418.900383
$ws.Range("I6").Value = 0.1723688299390553
$ws.Range("J6").Value = 0.1847819285778738
$ws.Range("M6").Value = 13.62041
$ws.Range("N6").Value = 40.86123000000001
$ws.Range("O6").Value = 0.1369489801107134
$ws.Range("P6").Value = 0.1486879047638899
$ws.Range("Q6").Value = 1901.86498853901
$ws.Range("R6").Value = 17116.7848968511
$ws.Range("S6").Value = 0.02360573546303063
$ws.Range("T6").Value = 0.0274748377984748
$ws.Range("I7").Value = 0.1574979484290947
$ws.Range("J7").Value = 0.1688401242154775
$ws.Range("M7").Value = 38.10639333333333
$ws.Range("N7").Value = 114.31918
$ws.Range("O7").Value = 0.3831479157160237
$ws.Range("P7").Value = 0.4159903984418967
$ws.Range("Q7").Value = 4861.871946600681
$ws.Range("R7").Value = 43756.84751940612
$ws.Range("S7").Value = 0.06034501067015744
$ws.Range("T7").Value = 0.07023587054537582
$ws.Range("I8").Value = 0.1574979484290947
$ws.Range("J8").Value = 0.1688401242154775
$ws.Range("O8").Value = 0.09199521176963764
$ws.Range("P8").Value = 0.09988081163714851
$ws.Range("S8").Value = 0.01448905711901804
$ws.Range("T8").Value = 0.01686388864355887
$ws.Range("I9").Value = 0.1574979484290947
$ws.Range("J9").Value = 0.1688401242154775
$ws.Range("M9").Value = 15.023598
$ws.Range("N9").Value = 45.070794
$ws.Range("O9").Value = 0.1510575983904562
$ws.Range("P9").Value = 0.1640058785774412
$ws.Range("Q9").Value = 1916.812462787244
$ws.Range("R9").Value = 17251.3121650852
$ws.Range("S9").Value = 0.02379126184112297
$ws.Range("T9").Value = 0.0276907729110837
$ws.Range("I10").Value = 0.1574979484290947
$ws.Range("J10").Value = 0.1688401242154775
$ws.Range("M10").Value = 23.556204
$ws.Range("N10").Value = 47.112408
$ws.Range("O10").Value = 0.236850294013169
$ws.Range("P10").Value = 0.1714350065796238
$ws.Range("Q10").Value = 3005.460170270713
$ws.Range("R10").Value = 18032.76102162428
$ws.Range("S10").Value = 0.03730343539190201
$ws.Range("T10").Value = 0.0289451078057849
$ws.Range("I11").Value = 0.1574979484290947
$ws.Range("J11").Value = 0.1688401242154775
$ws.Range("M11").Value = 13.62041
$ws.Range("N11").Value = 40.86123000000001
$ws.Range("O11").Value = 0.1369489801107134
$ws.Range("P11").Value = 0.1486879047638899
$ws.Range("Q11").Value = 1737.784226938981
$ws.Range("R11").Value = 15640.05804245082
$ws.Range("S11").Value = 0.02156918340689426
$ws.Range("T11").Value = 0.02510448430967426
$ws.Range("G12").Value = 180.251373
$ws.Range("H12").Value = 540.7541189999999
$ws.Range("I12").Value = 0.2225091180610228
$ws.Range("J12").Value = 0.2385330571427265
$ws.Range("M12").Value = 38.10639333333333
$ws.Range("N12").Value = 114.31918
$ws.Range("O12").Value = 0.3831479157160237
$ws.Range("P12").Value = 0.4159903984418967
$ws.Range("Q12").Value = 6868.729718411379
$ws.Range("R12").Value = 61818.56746570241
$ws.Range("S12").Value = 0.08525390481289154
$ws.Range("T12").Value = 0.09922746148236651
$ws.Range("G13").Value = 180.251373
$ws.Range("H13").Value = 540.7541189999999
$ws.Range("I13").Value = 0.2225091180610228
$ws.Range("J13").Value = 0.2385330571427265
$ws.Range("O13").Value = 0.09199521176963764
$ws.Range("P13").Value = 0.09988081163714851
$ws.Range("Q13").Value = 1649.207053241532
$ws.Range("R13").Value = 14842.86347917379
$ws.Range("S13").Value = 0.0204697734366991
$ws.Range("T13").Value = 0.02382487534970585
$ws.Range("G14").Value = 180.251373
$ws.Range("H14").Value = 540.7541189999999
$ws.Range("I14").Value = 0.2225091180610228
$ws.Range("J14").Value = 0.2385330571427265
$ws.Range("M14").Value = 15.023598
$ws.Range("N14").Value = 45.070794
$ws.Range("O14").Value = 0.1510575983904562
$ws.Range("P14").Value = 0.1640058785774412
$ws.Range("Q14").Value = 2708.024166900053
$ws.Range("R14").Value = 24372.21750210048
$ws.Range("S14").Value = 0.03361169299427658
$ws.Range("T14").Value = 0.03912082360645586
$ws.Range("G15").Value = 180.251373
$ws.Range("H15").Value = 540.7541189999999
$ws.Range("I15").Value = 0.2225091180610228
$ws.Range("J15").Value = 0.2385330571427265
$ws.Range("M15").Value = 23.556204
$ws.Range("N15").Value = 47.112408
$ws.Range("O15").Value = 0.236850294013169
$ws.Range("P15").Value = 0.1714350065796238
$ws.Range("Q15").Value = 4246.038113668092
$ws.Range("R15").Value = 25476.22868200855
$ws.Range("S15").Value = 0.05270135003336418
$ws.Range("T15").Value = 0.04089291622072111
$ws.Range("G16").Value = 180.251373
$ws.Range("H16").Value = 540.7541189999999
$ws.Range("I16").Value = 0.2225091180610228
$ws.Range("J16").Value = 0.2385330571427265
$ws.Range("M16").Value = 13.62041
$ws.Range("N16").Value = 40.86123000000001
$ws.Range("O16").Value = 0.1369489801107134
$ws.Range("P16").Value = 0.1486879047638899
$ws.Range("Q16").Value = 2455.09760332293
$ws.Range("R16").Value = 22095.87842990637
$ws.Range("S16").Value = 0.03047239678379139
$ws.Range("T16").Value = 0.03546698048347723
$ws.Range("G17").Value = 163.257347
$ws.Range("H17").Value = 326.514694
$ws.Range("I17").Value = 0.2015310490752954
$ws.Range("J17").Value = 0.1440295051397322
$ws.Range("M17").Value = 38.10639333333333
$ws.Range("N17").Value = 114.31918
$ws.Range("O17").Value = 0.3831479157160237
$ws.Range("P17").Value = 0.4159903984418967
$ws.Range("Q17").Value = 6221.148679338487
$ws.Range("R17").Value = 37326.89207603092
$ws.Range("S17").Value = 0.07721620140526313
$ws.Range("T17").Value = 0.05991489123046641
$ws.Range("G18").Value = 163.257347
$ws.Range("H18").Value = 326.514694
$ws.Range("I18").Value = 0.2015310490752954
$ws.Range("J18").Value = 0.1440295051397322
$ws.Range("O18").Value = 0.09199521176963764
$ws.Range("P18").Value = 0.09988081163714851
$ws.Range("Q18").Value = 1493.720484258948
$ws.Range("R18").Value = 8962.322905553689
$ws.Range("S18").Value = 0.01853989153783904
$ws.Range("T18").Value = 0.01438578387305331
$ws.Range("G19").Value = 163.257347
$ws.Range("H19").Value = 326.514694
$ws.Range("I19").Value = 0.2015310490752954
$ws.Range("J19").Value = 0.1440295051397322
$ws.Range("M19").Value = 15.023598
$ws.Range("N19").Value = 45.070794
$ws.Range("O19").Value = 0.1510575983904562
$ws.Range("P19").Value = 0.1640058785774412
$ws.Range("Q19").Value = 2452.712751874506
$ws.Range("R19").Value = 14716.27651124704
$ws.Range("S19").Value = 0.03044279627442329
$ws.Range("T19").Value = 0.02362168553151587
$ws.Range("G20").Value = 163.257347
$ws.Range("H20").Value = 326.514694
$ws.Range("I20").Value = 0.2015310490752954
$ws.Range("J20").Value = 0.1440295051397322
$ws.Range("M20").Value = 23.556204
$ws.Range("N20").Value = 47.112408
$ws.Range("O20").Value = 0.236850294013169
$ws.Range("P20").Value = 0.1714350065796238
$ws.Range("Q20").Value = 3845.723370430788
$ws.Range("R20").Value = 15382.89348172315
$ws.Range("S20").Value = 0.04773268822626611
$ws.Range("T20").Value = 0.02469169916128996
$ws.Range("G21").Value = 163.257347
$ws.Range("H21").Value = 326.514694
$ws.Range("I21").Value = 0.2015310490752954
$ws.Range("J21").Value = 0.1440295051397322
$ws.Range("M21").Value = 13.62041
$ws.Range("N21").Value = 40.86123000000001
$ws.Range("O21").Value = 0.1369489801107134
$ws.Range("P21").Value = 0.1486879047638899
$ws.Range("Q21").Value = 2223.63200165227
$ws.Range("R21").Value = 13341.79200991362
$ws.Range("S21").Value = 0.02759947163150383
$ws.Range("T21").Value = 0.02141544534340669
$ws.Range("G22").Value = 199.356374
$ws.Range("H22").Value = 598.069122
$ws.Range("I22").Value = 0.2460930544955318
$ws.Range("J22").Value = 0.2638153849241901
$ws.Range("M22").Value = 38.10639333333333
$ws.Range("N22").Value = 114.31918
$ws.Range("O22").Value = 0.3831479157160237
$ws.Range("P22").Value = 0.4159903984418967
$ws.Range("Q22").Value = 7596.752401151106
$ws.Range("R22").Value = 68370.77161035995
$ws.Range("S22").Value = 0.09429004090215284
$ws.Range("T22").Value = 0.1097446670897162
$ws.Range("G23").Value = 199.356374
$ws.Range("H23").Value = 598.069122
$ws.Range("I23").Value = 0.2460930544955318
$ws.Range("J23").Value = 0.2638153849241901
$ws.Range("O23").Value = 0.09199521176963764
$ws.Range("P23").Value = 0.09988081163714851
$ws.Range("Q23").Value = 1824.007954211016
$ws.Range("R23").Value = 16416.07158789914
$ws.Range("S23").Value = 0.02263938266335342
$ws.Range("T23").Value = 0.02635009476859486
$ws.Range("G24").Value = 199.356374
$ws.Range("H24").Value = 598.069122
$ws.Range("I24").Value = 0.2460930544955318
$ws.Range("J24").Value = 0.2638153849241901
$ws.Range("M24").Value = 15.023598
$ws.Range("N24").Value = 45.070794
$ws.Range("O24").Value = 0.1510575983904562
$ws.Range("P24").Value = 0.1640058785774412
$ws.Range("Q24").Value = 2995.050021713652
$ws.Range("R24").Value = 26955.45019542287
$ws.Range("S24").Value = 0.03717422579266669
$ws.Range("T24").Value = 0.04326727398673764
$ws.Range("G25").Value = 199.356374
$ws.Range("H25").Value = 598.069122
$ws.Range("I25").Value = 0.2460930544955318
$ws.Range("J25").Value = 0.2638153849241901
$ws.Range("M25").Value = 23.556204
$ws.Range("N25").Value = 47.112408
$ws.Range("O25").Value = 0.236850294013169
$ws.Range("P25").Value = 0.1714350065796238
$ws.Range("Q25").Value = 4696.079414644296
$ws.Range("R25").Value = 28176.47648786578
$ws.Range("S25").Value = 0.05828721231186552
$ws.Range("T25").Value = 0.04522719225028453
$ws.Range("G26").Value = 199.356374
$ws.Range("H26").Value = 598.069122
$ws.Range("I26").Value = 0.2460930544955318
$ws.Range("J26").Value = 0.2638153849241901
$ws.Range("M26").Value = 13.62041
$ws.Range("N26").Value = 40.86123000000001
$ws.Range("O26").Value = 0.1369489801107134
$ws.Range("P26").Value = 0.1486879047638899
$ws.Range("Q26").Value = 2715.31554999334
$ws.Range("R26").Value = 24437.83994994006
$ws.Range("S26").Value = 0.03370219282549328
$ws.Range("T26").Value = 0.03922615682885693
